# Converts an "RRGGBB" hex string into the little-endian BGR integer that
# the PowerPoint COM object model uses for RGB color values (same encoding
# as VBA's RGB() function / ThemeColor.RGB / Long color values).
function HexToComColor([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# --- 1. Table on slide 16 switches to a different built-in table style ---
$slide = $p.Slides.Item(16)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{81150F83-92B4-4DDF-8D75-8C8E18E06030}")
    }
}

# --- 2. Presentation theme switches from "Integral" to "Office Theme" ---
# Target ("Office Theme") color scheme values.
$officeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$scheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $scheme.Count; $i++) {
    $scheme.Item($i).RGB = HexToComColor $officeColors[$i - 1]
}
